$wb = $excel.ActiveWorkbook
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsElem = $wb.Worksheets.Item("Elements")

# ------------------------------------------------------------------
# 1. Update the IG url (pythia -> cicada) and regeneration date on the
#    Metadata sheet.
# ------------------------------------------------------------------
$wsMeta.Range("B2").Value = "http://fhirfli.dev/fhir/ig/cicada/StructureDefinition/preferred-interval-reason"
$wsMeta.Range("B8").Value = "2026-02-11T14:37:07-05:00"

# ------------------------------------------------------------------
# 2. Insert a new "Jurisdiction" property row right after "Contact"
#    (row 11), pushing every following row down by one. We shift the
#    values manually (instead of Rows.Insert) so the existing cell
#    style ("s=2") is preserved on every moved row instead of the
#    engine minting a brand new default style for an inserted row.
# ------------------------------------------------------------------
for ($r = 20; $r -ge 11; $r--) {
    $dst = $r + 1
    $aVal = $wsMeta.Cells.Item($r, 1).Value2
    $bVal = $wsMeta.Cells.Item($r, 2).Value2

    if ($aVal -eq $null) {
        $wsMeta.Cells.Item($dst, 1).Value = ""
    } elseif ($aVal -eq "false" -or $aVal -eq "true" -or $aVal -eq "0" -or $aVal -eq "1") {
        # Force literal text instead of letting Excel coerce to a
        # Boolean/number (matches the original shared-string type).
        $wsMeta.Cells.Item($dst, 1).Formula = "=T(""" + $aVal + """)"
        $wsMeta.Cells.Item($dst, 1).Copy()
        $wsMeta.Cells.Item($dst, 1).PasteSpecial(-4163)
    } else {
        $wsMeta.Cells.Item($dst, 1).Value = $aVal
    }

    if ($bVal -eq $null) {
        $wsMeta.Cells.Item($dst, 2).Value = ""
    } elseif ($bVal -eq "false" -or $bVal -eq "true" -or $bVal -eq "0" -or $bVal -eq "1") {
        $wsMeta.Cells.Item($dst, 2).Formula = "=T(""" + $bVal + """)"
        $wsMeta.Cells.Item($dst, 2).Copy()
        $wsMeta.Cells.Item($dst, 2).PasteSpecial(-4163)
    } else {
        $wsMeta.Cells.Item($dst, 2).Value = $bVal
    }
}

$wsMeta.Range("A11").Value = "Jurisdiction"
$wsMeta.Range("B11").Value = ""

# ------------------------------------------------------------------
# 3. Update the two other pythia -> cicada URLs on the Elements sheet
#    (Extension.url fixed value + the interval-reason value set url).
# ------------------------------------------------------------------
$wsElem.Range("R5").Value = "http://fhirfli.dev/fhir/ig/cicada/StructureDefinition/preferred-interval-reason"
$wsElem.Range("Z6").Value = "http://fhirfli.dev/fhir/ig/cicada/ValueSet/interval-reason"
